$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 2.098888
$ws.Range("H2").Value = 6.296664
$ws.Range("I2").Value = 0.1082453658858517
$ws.Range("J2").Value = 0.1082453658858517
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.8369776666666665
$ws.Range("N2").Value = 2.510933
$ws.Range("O2").Value = 0.0694586718035551
$ws.Range("P2").Value = 0.06945867180355511
$ws.Range("Q2").Value = 1.756722380834666
$ws.Range("R2").Value = 15.810501427512
$ws.Range("S2").Value = 0.007518579343321111
$ws.Range("T2").Value = 0.00751857934332111
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 2.098888
$ws.Range("H3").Value = 6.296664
$ws.Range("I3").Value = 0.1082453658858517
$ws.Range("J3").Value = 0.1082453658858517
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 2.427350333333333
$ws.Range("N3").Value = 7.282051
$ws.Range("O3").Value = 0.2014397000898671
$ws.Range("P3").Value = 0.2014397000898671
$ws.Range("Q3").Value = 5.094736486429334
$ws.Range("R3").Value = 45.852628377864
$ws.Range("S3").Value = 0.0218049140401639
$ws.Range("T3").Value = 0.02180491404016389
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 2.098888
$ws.Range("H4").Value = 6.296664
$ws.Range("I4").Value = 0.1082453658858517
$ws.Range("J4").Value = 0.1082453658858517
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 8.785681666666667
$ws.Range("N4").Value = 26.357045
$ws.Range("O4").Value = 0.7291016281065776
$ws.Range("P4").Value = 0.7291016281065776
$ws.Range("Q4").Value = 18.44016182198667
$ws.Range("R4").Value = 165.96145639788
$ws.Range("S4").Value = 0.07892187250236665
$ws.Range("T4").Value = 0.07892187250236664
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 15.87514366666667
$ws.Range("H5").Value = 47.625431
$ws.Range("I5").Value = 0.8187243600843848
$ws.Range("J5").Value = 0.8187243600843847
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.8369776666666665
$ws.Range("N5").Value = 2.510933
$ws.Range("O5").Value = 0.0694586718035551
$ws.Range("P5").Value = 0.06945867180355511
$ws.Range("Q5").Value = 13.28714070412478
$ws.Range("R5").Value = 119.584266337123
$ws.Range("S5").Value = 0.05686750662467695
$ws.Range("T5").Value = 0.05686750662467695
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 15.87514366666667
$ws.Range("H6").Value = 47.625431
$ws.Range("I6").Value = 0.8187243600843848
$ws.Range("J6").Value = 0.8187243600843847
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 2.427350333333333
$ws.Range("N6").Value = 7.282051
$ws.Range("O6").Value = 0.2014397000898671
$ws.Range("P6").Value = 0.2014397000898671
$ws.Range("Q6").Value = 38.53453527099789
$ws.Range("R6").Value = 346.810817438981
$ws.Range("S6").Value = 0.1649235895516669
$ws.Range("T6").Value = 0.1649235895516668
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 15.87514366666667
$ws.Range("H7").Value = 47.625431
$ws.Range("I7").Value = 0.8187243600843848
$ws.Range("J7").Value = 0.8187243600843847
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 8.785681666666667
$ws.Range("N7").Value = 26.357045
$ws.Range("O7").Value = 0.7291016281065776
$ws.Range("P7").Value = 0.7291016281065776
$ws.Range("Q7").Value = 139.4739586679328
$ws.Range("R7").Value = 1255.265628011395
$ws.Range("S7").Value = 0.5969332639080408
$ws.Range("T7").Value = 0.5969332639080408
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 1.416064
$ws.Range("H8").Value = 4.248192
$ws.Range("I8").Value = 0.07303027402976368
$ws.Range("J8").Value = 0.07303027402976367
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.8369776666666665
$ws.Range("N8").Value = 2.510933
$ws.Range("O8").Value = 0.0694586718035551
$ws.Range("P8").Value = 0.06945867180355511
$ws.Range("Q8").Value = 1.185213942570666
$ws.Range("R8").Value = 10.666925483136
$ws.Range("S8").Value = 0.005072585835557049
$ws.Range("T8").Value = 0.005072585835557049
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 1.416064
$ws.Range("H9").Value = 4.248192
$ws.Range("I9").Value = 0.07303027402976368
$ws.Range("J9").Value = 0.07303027402976367
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 2.427350333333333
$ws.Range("N9").Value = 7.282051
$ws.Range("O9").Value = 0.2014397000898671
$ws.Range("P9").Value = 0.2014397000898671
$ws.Range("Q9").Value = 3.437283422421333
$ws.Range("R9").Value = 30.935550801792
$ws.Range("S9").Value = 0.01471119649803641
$ws.Range("T9").Value = 0.01471119649803641
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 1.416064
$ws.Range("H10").Value = 4.248192
$ws.Range("I10").Value = 0.07303027402976368
$ws.Range("J10").Value = 0.07303027402976367
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 8.785681666666667
$ws.Range("N10").Value = 26.357045
$ws.Range("O10").Value = 0.7291016281065776
$ws.Range("P10").Value = 0.7291016281065776
$ws.Range("Q10").Value = 12.44108752362667
$ws.Range("R10").Value = 111.96978771264
$ws.Range("S10").Value = 0.05324649169617022
$ws.Range("T10").Value = 0.0532464916961702
